$d = $word.ActiveDocument

function Insert-PlainRun {
    param($targetCollapsedRange, $text)
    $start = $targetCollapsedRange.Start
    $targetCollapsedRange.InsertAfter($text)
    $newEnd = $start + $text.Length
    $result = $d.Range($newEnd, $newEnd)
    return $result
}

function Paste-FormattedRun {
    param($sourceRange, $targetCollapsedRange, $text)
    $start = $targetCollapsedRange.Start
    $sourceRange.Copy()
    $targetCollapsedRange.Paste()
    $srcLen = $sourceRange.Text.Length
    $newEndTmp = $start + $srcLen
    $newRange = $d.Range($start, $newEndTmp)
    $newRange.Text = $text
    $newEnd = $start + $text.Length
    $result = $d.Range($newEnd, $newEnd)
    return $result
}

# The existing "_GoBack" bookmark sits right after " lop -p" in the "For details
# log" paragraph. We delete it now and re-create it after all the new content has
# been inserted, at its new final location (end of the "For rename file" line).
if ($d.Bookmarks.Exists("_GoBack")) {
    $gb = $d.Bookmarks.Item("_GoBack")
    $gb.Delete()
}

# Locate the paragraph "For details log: git lop -p" (paragraph 18).
$p18 = $d.Paragraphs.Item(18)

# Template range carrying the run formatting used throughout this section:
# Times New Roman / Bold / 24 (12pt) / Complex-script 24 / lang en-US.
$template = $p18.Range.Duplicate
$template.Find.Execute("git") | Out-Null

# 1) " lop -p"  ->  " lop " + "\u2013" + "p"   (3 runs, same formatting)
$r = $p18.Range.Duplicate
$r.Find.Execute(" lop -p") | Out-Null
$r.Text = " lop "
$r.Collapse(0)

$dashChar = [char]0x2013
$dashStr = [string]$dashChar
$r = Paste-FormattedRun $template $r $dashStr
$r = Paste-FormattedRun $template $r "p"

# 2) Insert five new paragraphs after the "For details log" paragraph.
$p18 = $d.Paragraphs.Item(18)
$endR = $d.Range($p18.Range.End - 1, $p18.Range.End - 1)
$endR.InsertParagraphAfter()

$p19 = $d.Paragraphs.Item(19)
$ins = $d.Range($p19.Range.End - 1, $p19.Range.End - 1)
$ins = Insert-PlainRun $ins "For queit use: Q"

$endR = $d.Range($p19.Range.End - 1, $p19.Range.End - 1)
$endR.InsertParagraphAfter()

$p20 = $d.Paragraphs.Item(20)
$ins = $d.Range($p20.Range.End - 1, $p20.Range.End - 1)
$ins = Insert-PlainRun $ins "For cronocila order: gitreflog"

$endR = $d.Range($p20.Range.End - 1, $p20.Range.End - 1)
$endR.InsertParagraphAfter()

$p21 = $d.Paragraphs.Item(21)
$ins = $d.Range($p21.Range.End - 1, $p21.Range.End - 1)
$ins = Insert-PlainRun $ins "For Git specific hash details: git show <Hashcode>"

$endR = $d.Range($p21.Range.End - 1, $p21.Range.End - 1)
$endR.InsertParagraphAfter()

$p22 = $d.Paragraphs.Item(22)
$ins = $d.Range($p22.Range.End - 1, $p22.Range.End - 1)
$ins = Insert-PlainRun $ins "For delete the file: git rm <file name>"

$endR = $d.Range($p22.Range.End - 1, $p22.Range.End - 1)
$endR.InsertParagraphAfter()

$p23 = $d.Paragraphs.Item(23)
$ins = $d.Range($p23.Range.End - 1, $p23.Range.End - 1)
$ins = Insert-PlainRun $ins "For rename file: git mv <old file name with extension>  <new file name with extension>"

# 3) Three blank paragraphs after the "For rename file" paragraph.
$endR = $d.Range($p23.Range.End - 1, $p23.Range.End - 1)
$endR.InsertParagraphAfter()
$p24 = $d.Paragraphs.Item(24)

$endR = $d.Range($p24.Range.End - 1, $p24.Range.End - 1)
$endR.InsertParagraphAfter()
$p25 = $d.Paragraphs.Item(25)

$endR = $d.Range($p25.Range.End - 1, $p25.Range.End - 1)
$endR.InsertParagraphAfter()

# 4) Re-create the "_GoBack" bookmark at the end of the "For rename file" paragraph.
$p23 = $d.Paragraphs.Item(23)
$bmPos = $p23.Range.End - 1
$bmRange = $d.Range($bmPos, $bmPos)
$d.Bookmarks.Add("_GoBack", $bmRange)

# 5) Mark the (second) "For set user name:" paragraph (the user.email one) with a
# lastRenderedPageBreak right before its first run, as in the target revision.
$all = $d.Content
$pageBreakTarget = $all.Duplicate
$pageBreakTarget.Find.MatchCase = $true
$pageBreakTarget.Find.Execute("For set user name:  Git config --global user.email") | Out-Null
$pageBreakTarget.Collapse(1)
$pageBreakTarget.InsertParagraphBefore()
